$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2599.963
$ws.Range("I51").Value = 2491.6667
$ws.Range("J51").Value = 3466.3333
$ws.Range("K51").Value = 2491.6667
$ws.Range("L51").Value = 3466.3333
$ws.Range("M51").Value = -2007.6667
$ws.Range("N51").Value = -4434.3333

$ws.Range("H58").Value = 330.46155
$ws.Range("I58").Value = 290.54544
$ws.Range("J58").Value = 550
$ws.Range("K58").Value = 871.63632
$ws.Range("L58").Value = 1650
$ws.Range("M58").Value = -721.63632
$ws.Range("N58").Value = -1950

$ws.Range("H62").Value = 12088.81
$ws.Range("I62").Value = 12468.5
$ws.Range("K62").Value = 12468.5
$ws.Range("M62").Value = -11844.5

$ws.Range("H65").Value = 12088.81
$ws.Range("I65").Value = 12468.5
$ws.Range("K65").Value = 62342.5
$ws.Range("M65").Value = -59222.5

$ws.Range("H88").Value = 8750
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()

$ws.Range("H91").Value = 8750
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()

$ws.Range("H116").Value = 12579.546
$ws.Range("I116").Value = 15876
$ws.Range("J116").Value = 9832.5
$ws.Range("K116").Value = 15876
$ws.Range("L116").Value = 9832.5
$ws.Range("M116").Value = -12434
$ws.Range("N116").Value = -16716.5

$ws.Range("H132").Value = 8359463.5
$ws.Range("I132").Value = 9010119
$ws.Range("J132").Value = 334716.66
$ws.Range("K132").Value = 27030357
$ws.Range("L132").Value = 1004149.98
$ws.Range("M132").Value = -27027827
$ws.Range("N132").Value = -1009209.98

$ws.Range("H137").Value = 1603.8422
$ws.Range("I137").Value = 1408.6
$ws.Range("K137").Value = 4225.799999999999
$ws.Range("M137").Value = -1675.799999999999

$ws.Range("H138").Value = 2926.2388
$ws.Range("I138").Value = 1274.303
$ws.Range("J138").Value = 4529.5884
$ws.Range("K138").Value = 3822.909000000001
$ws.Range("L138").Value = 13588.7652
$ws.Range("M138").Value = 1317.090999999999
$ws.Range("N138").Value = -23868.7652

$ws.Range("H141").Value = 2502
$ws.Range("I141").Value = 1903.1111
$ws.Range("K141").Value = 5709.3333
$ws.Range("M141").Value = -529.3333000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 9035.200000000001
$ws.Range("J45").Value = 5000
$ws.Range("L45").Value = 5000
$ws.Range("N45").Value = -5754

$ws.Range("H61").Value = 3830.9062
$ws.Range("I61").Value = 3516.0833
$ws.Range("J61").Value = 4775.375
$ws.Range("K61").Value = 3516.0833
$ws.Range("L61").Value = 4775.375
$ws.Range("M61").Value = -3304.0833
$ws.Range("N61").Value = -5199.375

$ws.Range("H122").Value = 1503.1052
$ws.Range("I122").Value = 1531.1111
$ws.Range("K122").Value = 4593.3333
$ws.Range("M122").Value = -2143.3333

$ws.Range("H124").Value = 14750
$ws.Range("J124").Value = 14750
$ws.Range("L124").Value = 14750
$ws.Range("N124").Value = -24570

$ws.Range("H132").Value = 3463.125
$ws.Range("I132").Value = 1885.8667
$ws.Range("J132").Value = 6091.8887
$ws.Range("K132").Value = 5657.6001
$ws.Range("L132").Value = 18275.6661
$ws.Range("M132").Value = -3127.6001
$ws.Range("N132").Value = -23335.6661

$ws.Range("H133").Value = 147999.67
$ws.Range("J133").Value = 211499.5
$ws.Range("L133").Value = 211499.5
$ws.Range("N133").Value = -216559.5

$ws.Range("H136").Value = 3830.9062
$ws.Range("I136").Value = 3516.0833
$ws.Range("J136").Value = 4775.375
$ws.Range("K136").Value = 10548.2499
$ws.Range("L136").Value = 14326.125
$ws.Range("M136").Value = -7998.249899999999
$ws.Range("N136").Value = -19426.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 27329
$ws.Range("I82").Value = 4743.5
$ws.Range("J82").Value = 72500
$ws.Range("K82").Value = 4743.5
$ws.Range("L82").Value = 72500
$ws.Range("M82").Value = -4360.5
$ws.Range("N82").Value = -73266

$ws.Range("H85").Value = 27329
$ws.Range("I85").Value = 4743.5
$ws.Range("J85").Value = 72500
$ws.Range("K85").Value = 4743.5
$ws.Range("L85").Value = 72500
$ws.Range("M85").Value = -3417.5
$ws.Range("N85").Value = -75152

$ws.Range("H95").Value = 15331.143
$ws.Range("J95").Value = 15331.143
$ws.Range("L95").Value = 15331.143
$ws.Range("N95").Value = -20823.143

$ws.Range("H134").Value = 2235.7568
$ws.Range("I134").Value = 1957.6897
$ws.Range("K134").Value = 5873.0691
$ws.Range("M134").Value = -3338.0691

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 12554
$ws.Range("J28").Value = 12748.375
$ws.Range("L28").Value = 12748.375
$ws.Range("N28").Value = -13238.375

$ws.Range("H58").Value = 2414.4595
$ws.Range("I58").Value = 2501.1562
$ws.Range("J58").Value = 1859.6
$ws.Range("K58").Value = 2501.1562
$ws.Range("L58").Value = 1859.6
$ws.Range("M58").Value = -2298.1562
$ws.Range("N58").Value = -2265.6

$ws.Range("H132").Value = 2934.3
$ws.Range("I132").Value = 3072.0605
$ws.Range("J132").Value = 2284.8572
$ws.Range("K132").Value = 9216.181500000001
$ws.Range("L132").Value = 6854.571599999999
$ws.Range("M132").Value = -6686.181500000001
$ws.Range("N132").Value = -11914.5716

$ws.Range("H136").Value = 2414.4595
$ws.Range("I136").Value = 2501.1562
$ws.Range("J136").Value = 1859.6
$ws.Range("K136").Value = 7503.4686
$ws.Range("L136").Value = 5578.799999999999
$ws.Range("M136").Value = -4953.4686
$ws.Range("N136").Value = -10678.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 3960.1428
$ws.Range("I81").Value = 941
$ws.Range("J81").Value = 6224.5
$ws.Range("K81").Value = 2823
$ws.Range("L81").Value = 18673.5
$ws.Range("M81").Value = -1700
$ws.Range("N81").Value = -20919.5

$ws.Range("H84").Value = 3960.1428
$ws.Range("I84").Value = 941
$ws.Range("J84").Value = 6224.5
$ws.Range("K84").Value = 8469
$ws.Range("L84").Value = 56020.5
$ws.Range("M84").Value = -2853
$ws.Range("N84").Value = -67252.5

$ws.Range("H93").Value = 6500
$ws.Range("I93").Value = 5000
$ws.Range("J93").Value = 8000
$ws.Range("K93").Value = 15000
$ws.Range("L93").Value = 24000
$ws.Range("M93").Value = -13128
$ws.Range("N93").Value = -27744

$ws.Range("H103").Value = 5881.6665
$ws.Range("I103").Value = 150
$ws.Range("J103").Value = 7028
$ws.Range("K103").Value = 450
$ws.Range("L103").Value = 21084
$ws.Range("M103").Value = 429
$ws.Range("N103").Value = -22842

$ws.Range("H113").Value = 537.9474
$ws.Range("J113").Value = 545.6111
$ws.Range("L113").Value = 1636.8333
$ws.Range("N113").Value = -5976.8333

$ws.Range("H121").Value = 219.8
$ws.Range("I121").Value = 149.75
$ws.Range("J121").Value = 500
$ws.Range("K121").Value = 449.25
$ws.Range("L121").Value = 1500
$ws.Range("M121").Value = 860.75
$ws.Range("N121").Value = -4120

$ws.Range("H129").Value = 1161.5
$ws.Range("I129").Value = 907.2
$ws.Range("J129").Value = 2433
$ws.Range("K129").Value = 2721.6
$ws.Range("L129").Value = 7299
$ws.Range("M129").Value = 2278.4
$ws.Range("N129").Value = -17299

$ws.Range("H131").Value = 39282.85
$ws.Range("I131").Value = 111988.22
$ws.Range("J131").Value = 2930.1667
$ws.Range("K131").Value = 335964.66
$ws.Range("L131").Value = 8790.500100000001
$ws.Range("M131").Value = -330924.66
$ws.Range("N131").Value = -18870.5001

$ws.Range("H138").Value = 21747948
$ws.Range("I138").Value = 45461564
$ws.Range("K138").Value = 136384692
$ws.Range("M138").Value = -136379552

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H94").Value = 100000
$ws.Range("J94").Value = 100000
$ws.Range("L94").Value = 100000
$ws.Range("N94").Value = -101352

$ws.Range("H102").Value = 1065.1818
$ws.Range("I102").Value = 980.1053000000001
$ws.Range("J102").Value = 1604
$ws.Range("K102").Value = 980.1053000000001
$ws.Range("L102").Value = 1604
$ws.Range("M102").Value = 641.8946999999999
$ws.Range("N102").Value = -4848

$ws.Range("H126").Value = 12310.728
$ws.Range("J126").Value = 3777.5
$ws.Range("L126").Value = 11332.5
$ws.Range("N126").Value = -16272.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3837.6191
$ws.Range("J61").Value = 2699.5
$ws.Range("L61").Value = 2699.5
$ws.Range("N61").Value = -3103.5

$ws.Range("H113").Value = 3837.6191
$ws.Range("J113").Value = 2699.5
$ws.Range("L113").Value = 2699.5
$ws.Range("N113").Value = -7039.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7675.3477
$ws.Range("I81").Value = 12172.223
$ws.Range("J81").Value = 4784.5
$ws.Range("K81").Value = 24344.446
$ws.Range("L81").Value = 9569
$ws.Range("M81").Value = -23283.446
$ws.Range("N81").Value = -11691

$ws.Range("H84").Value = 7675.3477
$ws.Range("I84").Value = 12172.223
$ws.Range("J84").Value = 4784.5
$ws.Range("K84").Value = 121722.23
$ws.Range("L84").Value = 47845
$ws.Range("M84").Value = -116418.23
$ws.Range("N84").Value = -58453

$ws.Range("H132").Value = 3706.7896
$ws.Range("I132").Value = 3652.8276
$ws.Range("K132").Value = 10958.4828
$ws.Range("M132").Value = -8428.4828

$ws.Range("H135").Value = 114999
$ws.Range("J135").Value = 114999
$ws.Range("L135").Value = 114999
$ws.Range("N135").Value = -125139
